# "Generate Report for Handback" - updates the localization-status workbook
# after a handback has completed: status text changes, handback datetime /
# handback file / target file columns get populated, and a hyperlink to the
# handed-back source file is added for each locale sheet.

$wb = $excel.ActiveWorkbook

$statusHandedBack = "Handed back: in sync with en-US"
$sourceFileName   = "29d5be5b-0675-4708-817e-4ea116377819.md"
$sourceFileUrl    = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bee926a3910001eeb98ee539c84eada27b0d4988/e2e/29d5be5b-0675-4708-817e-4ea116377819.md"

# Blue, underlined "hyperlink" look-and-feel used elsewhere in the workbook
# (RGB 0x64,0x95,0xED stored as a BGR long, the way VBA/COM colors work).
$hyperlinkColor = 15570276

function Set-HandbackWidth($ws, $col, $target) {
    # ColumnWidth is expressed in "characters" and gets snapped by Excel to
    # whole-pixel boundaries, so we dial in the closest value that rounds to
    # the desired stored width.
    $ws.Columns.Item($col).ColumnWidth = $target - 0.8333333333333334
}

# ---------------------------------------------------------------------
# Overview sheet: the per-locale status cells just reflect the new status
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $statusHandedBack
$overview.Range("F2").Value = $statusHandedBack
Set-HandbackWidth $overview 5 29.9777050018311
Set-HandbackWidth $overview 6 29.9777050018311

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $statusHandedBack
$zhcn.Range("L2").Value = "2017-02-09 16:22:00"

$zhcn.Range("J2").Value = $sourceFileName
$zhcn.Range("K2").Value = "29d5be5b-0675-4708-817e-4ea116377819.94e6e24afd5850ea8a21de4ede714b453f008a3f.zh-cn.xlf"

$zhcn.Hyperlinks.Add($zhcn.Range("J2"), $sourceFileUrl, "", "", $sourceFileName)
$zhcn.Range("J2").Font.Underline = $true
$zhcn.Range("J2").Font.Color = $hyperlinkColor

Set-HandbackWidth $zhcn 3 29.9777050018311
Set-HandbackWidth $zhcn 10 40
Set-HandbackWidth $zhcn 11 40

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $statusHandedBack
$dede.Range("L2").Value = "2017-02-09 16:22:25"

$dede.Range("J2").Value = $sourceFileName
$dede.Range("K2").Value = "29d5be5b-0675-4708-817e-4ea116377819.94e6e24afd5850ea8a21de4ede714b453f008a3f.de-de.xlf"

$dede.Hyperlinks.Add($dede.Range("J2"), $sourceFileUrl, "", "", $sourceFileName)
$dede.Range("J2").Font.Underline = $true
$dede.Range("J2").Font.Color = $hyperlinkColor

Set-HandbackWidth $dede 3 29.9777050018311
Set-HandbackWidth $dede 10 40
Set-HandbackWidth $dede 11 40
